$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")
$seriesSheet = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: append new weekly observation row 94 ---
$dataSheet.Range("A93").Copy($dataSheet.Range("A94"))
$dataSheet.Range("A94").Value = 45119
$dataSheet.Range("B94").Value = 3166.607

# --- SeriesInfo sheet: refresh metadata from latest FRED pull ---
# Helper: assign text values without letting Excel auto-coerce
# date-looking strings into date serials / changing the cell style.

$seriesSheet.Range("B3").NumberFormat = "@"
$seriesSheet.Range("B3").Value = "2023-07-20"
$seriesSheet.Range("B3").Style = "Normal"

$seriesSheet.Range("B4").NumberFormat = "@"
$seriesSheet.Range("B4").Value = "2023-07-20"
$seriesSheet.Range("B4").Style = "Normal"

$seriesSheet.Range("B7").NumberFormat = "@"
$seriesSheet.Range("B7").Value = "2023-07-12"
$seriesSheet.Range("B7").Style = "Normal"

$seriesSheet.Range("B14").NumberFormat = "@"
$seriesSheet.Range("B14").Value = "2023-07-13 15:35:27-05"
$seriesSheet.Range("B14").Style = "Normal"

$seriesSheet.Range("B15").Value = 77
